# Daily attendance processing - 2025-12-04 20:30:12
# Normalize the "Recorded By" (column G) entries so that the
# dnasr281@gmail.com account is listed last instead of first,
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$prefix = "dnasr281@gmail.com, "

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val.ToString().StartsWith($prefix)) {
        $rest = $val.ToString().Substring($prefix.Length)
        $cell.Value = $rest + ", dnasr281@gmail.com"
    }
}
